$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 110
$ws.Range("B110").Value = 6664682
$ws.Range("F110").Value = "Libertad Asuncion"
$ws.Range("G110").Value = "Olimpia Asuncion"
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = "H"
$ws.Range("K110").Value = 1.95
$ws.Range("L110").Value = 3.3
$ws.Range("M110").Value = 3.5
$ws.Range("N110").Value = 1.7
$ws.Range("O110").Value = 3.6
$ws.Range("P110").Value = 4.5
$ws.Range("Q110").Value = -0.75
$ws.Range("R110").Value = 1.9
$ws.Range("S110").Value = 1.9
$ws.Range("T110").Value = 2.5
$ws.Range("U110").Value = 1.85
$ws.Range("V110").Value = 1.95
$ws.Range("W110").Value = 0.7
$ws.Range("Y110").Value = -1
$ws.Range("Z110").Value = 0.45
$ws.Range("AA110").Value = -0.5
$ws.Range("AB110").Value = -1
$ws.Range("AC110").Value = 0.95

# Row 111
$ws.Range("B111").Value = 6664683
$ws.Range("F111").Value = "Cerro Porteno"
$ws.Range("G111").Value = "General Caballero JLM"
$ws.Range("I111").Value = 2
$ws.Range("J111").Value = "A"
$ws.Range("K111").Value = 1.363
$ws.Range("L111").Value = 4.5
$ws.Range("M111").Value = 7
$ws.Range("N111").Value = 1.45
$ws.Range("O111").Value = 4.2
$ws.Range("P111").Value = 6
$ws.Range("Q111").Value = -1.25
$ws.Range("R111").Value = 2
$ws.Range("S111").Value = 1.8
$ws.Range("T111").Value = 2.75
$ws.Range("U111").Value = 1.775
$ws.Range("V111").Value = 2.025
$ws.Range("W111").Value = -1
$ws.Range("Y111").Value = 5
$ws.Range("Z111").Value = -1
$ws.Range("AA111").Value = 0.8
$ws.Range("AB111").Value = 0.3875
$ws.Range("AC111").Value = -0.5

# Row 122
$ws.Range("B122").Value = 6720873
$ws.Range("F122").Value = "Sportivo Luqueno"
$ws.Range("G122").Value = "Sportivo Trinidense"
$ws.Range("H122").Value = 2
$ws.Range("J122").Value = "D"
$ws.Range("K122").Value = 2.625
$ws.Range("L122").Value = 3.1
$ws.Range("M122").Value = 2.5
$ws.Range("N122").Value = 2.3
$ws.Range("O122").Value = 3.1
$ws.Range("P122").Value = 2.9
$ws.Range("Q122").Value = -0.25
$ws.Range("R122").Value = 2.025
$ws.Range("S122").Value = 1.775
$ws.Range("U122").Value = 1.95
$ws.Range("V122").Value = 1.85
$ws.Range("X122").Value = 2.1
$ws.Range("Y122").Value = -1
$ws.Range("Z122").Value = -0.5
$ws.Range("AA122").Value = 0.3875
$ws.Range("AB122").Value = 0.95

# Row 123
$ws.Range("B123").Value = 6720844
$ws.Range("F123").Value = "Guarani Asuncion"
$ws.Range("G123").Value = "Olimpia Asuncion"
$ws.Range("H123").Value = 1
$ws.Range("J123").Value = "A"
$ws.Range("K123").Value = 2.45
$ws.Range("L123").Value = 3
$ws.Range("M123").Value = 2.75
$ws.Range("N123").Value = 4
$ws.Range("O123").Value = 3.2
$ws.Range("P123").Value = 1.85
$ws.Range("Q123").Value = 0.5
$ws.Range("R123").Value = 1.875
$ws.Range("S123").Value = 1.925
$ws.Range("U123").Value = 1.925
$ws.Range("V123").Value = 1.875
$ws.Range("X123").Value = -1
$ws.Range("Y123").Value = 0.8500000000000001
$ws.Range("Z123").Value = -1
$ws.Range("AA123").Value = 0.925
$ws.Range("AB123").Value = 0.925

# Row 260
$ws.Range("B260").Value = 7493310
$ws.Range("F260").Value = "Libertad Asuncion"
$ws.Range("G260").Value = "Tacuary"
$ws.Range("H260").Value = 1
$ws.Range("J260").Value = "A"
$ws.Range("K260").Value = 1.363
$ws.Range("L260").Value = 5
$ws.Range("M260").Value = 7
$ws.Range("N260").Value = 1.571
$ws.Range("O260").Value = 4.2
$ws.Range("P260").Value = 4.75
$ws.Range("Q260").Value = -0.75
$ws.Range("R260").Value = 1.8
$ws.Range("S260").Value = 2
$ws.Range("T260").Value = 2.75
$ws.Range("U260").Value = 1.8
$ws.Range("V260").Value = 2
$ws.Range("W260").Value = -1
$ws.Range("Y260").Value = 3.75
$ws.Range("Z260").Value = -1
$ws.Range("AA260").Value = 1
$ws.Range("AB260").Value = 0.4
$ws.Range("AC260").Value = -0.5

# Row 261
$ws.Range("B261").Value = 7493431
$ws.Range("F261").Value = "Sportivo Trinidense"
$ws.Range("G261").Value = "Guairena FC"
$ws.Range("H261").Value = 7
$ws.Range("J261").Value = "H"
$ws.Range("K261").Value = 2.05
$ws.Range("L261").Value = 3.3
$ws.Range("M261").Value = 3.3
$ws.Range("N261").Value = 2.6
$ws.Range("O261").Value = 3.1
$ws.Range("P261").Value = 2.6
$ws.Range("Q261").Value = 0
$ws.Range("R261").Value = 1.925
$ws.Range("S261").Value = 1.875
$ws.Range("T261").Value = 2.5
$ws.Range("U261").Value = 2
$ws.Range("V261").Value = 1.8
$ws.Range("W261").Value = 1.6
$ws.Range("Y261").Value = -1
$ws.Range("Z261").Value = 0.925
$ws.Range("AA261").Value = -1
$ws.Range("AB261").Value = 1
$ws.Range("AC261").Value = -1

# Row 263
$ws.Range("B263").Value = 7493433
$ws.Range("F263").Value = "Sportivo Luqueno"
$ws.Range("G263").Value = "Nacional Asuncion"
$ws.Range("H263").Value = 1
$ws.Range("I263").Value = 1
$ws.Range("J263").Value = "D"
$ws.Range("K263").Value = 2.75
$ws.Range("L263").Value = 3.2
$ws.Range("M263").Value = 2.4
$ws.Range("N263").Value = 2.75
$ws.Range("O263").Value = 3.1
$ws.Range("P263").Value = 2.45
$ws.Range("Q263").Value = 0.25
$ws.Range("R263").Value = 1.75
$ws.Range("S263").Value = 2.05
$ws.Range("T263").Value = 2.25
$ws.Range("U263").Value = 2
$ws.Range("V263").Value = 1.8
$ws.Range("W263").Value = -1
$ws.Range("X263").Value = 2.1
$ws.Range("Z263").Value = 0.375
$ws.Range("AA263").Value = -0.5
$ws.Range("AB263").Value = -0.5
$ws.Range("AC263").Value = 0.4

# Row 264
$ws.Range("B264").Value = 7493312
$ws.Range("F264").Value = "Cerro Porteno"
$ws.Range("G264").Value = "Guarani Asuncion"
$ws.Range("H264").Value = 4
$ws.Range("I264").Value = 0
$ws.Range("J264").Value = "H"
$ws.Range("K264").Value = 1.7
$ws.Range("L264").Value = 3.6
$ws.Range("M264").Value = 4.333
$ws.Range("N264").Value = 1.727
$ws.Range("O264").Value = 3.75
$ws.Range("P264").Value = 4.2
$ws.Range("Q264").Value = -0.5
$ws.Range("R264").Value = 1.8
$ws.Range("S264").Value = 2
$ws.Range("T264").Value = 2.75
$ws.Range("U264").Value = 1.875
$ws.Range("V264").Value = 1.925
$ws.Range("W264").Value = 0.7270000000000001
$ws.Range("X264").Value = -1
$ws.Range("Z264").Value = 0.8
$ws.Range("AA264").Value = -1
$ws.Range("AB264").Value = 0.875
$ws.Range("AC264").Value = -1

# Row 301
$ws.Range("R301").Value = 1.95
$ws.Range("S301").Value = 1.85

# Row 304
$ws.Range("N304").Value = 2.625
$ws.Range("P304").Value = 2.45
$ws.Range("Q304").Value = 0
$ws.Range("R304").Value = 1.95
$ws.Range("S304").Value = 1.85
$ws.Range("U304").Value = 1.925
$ws.Range("V304").Value = 1.875

Write-Host "Applied league base update for 21-02-2024"
